$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = '*maa://24633 (56.52), *maa://30515 (69.9), *maa://34787 (72.97), ***maa://20792 (11.93), maa://39402 (91.23), ***maa://29083 (27.78)'
$ws.Range("AF2").Value = 'maa://25251 (92.17), ***maa://21730 (25.33), ***maa://39501 (17.24), **maa://36675 (50.0)'
$ws.Range("H3").Value = 'maa://21247 (98.54), *maa://22748 (60.0)'
$ws.Range("X4").Value = '**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (88.24), ***maa://36683 (28.26)'
$ws.Range("A8").Value = '更新日期：2025.02.09 13:17:45'
$ws.Range("X8").Value = 'maa://21411 (95.86)'
$ws.Range("D10").Value = '***maa://25695 (18.82), ***maa://34206 (20.0), ***maa://39951 (14.0), ***maa://39243 (28.57), *maa://45271 (56.0)'
$ws.Range("D13").Value = 'maa://24999 (92.02), maa://36673 (93.15), maa://25001 (85.71)'
$ws.Range("AF13").Value = '**maa://22737 (33.33), maa://39883 (91.18), *maa://39885 (55.17)'
$ws.Range("D15").Value = '*maa://22743 (77.51), maa://22734 (84.03), *maa://30808 (64.18), **maa://36048 (44.07), maa://45058 (90.0)'
$ws.Range("AF15").Value = 'maa://21364 (81.17), *maa://36666 (79.25), *maa://22766 (68.38)'
$ws.Range("D16").Value = 'maa://21441 (96.4), maa://36679 (94.12), maa://37650 (97.14)'
$ws.Range("T16").Value = 'maa://22729 (94.9), *maa://28648 (69.12), maa://36674 (80.85)'
$ws.Range("AF16").Value = '*maa://23911 (65.09), maa://27755 (93.48)'
$ws.Range("D20").Value = 'maa://21432 (90.24), maa://25198 (93.46), *maa://20795 (51.16), maa://36680 (93.75)'
$ws.Range("L20").Value = 'maa://41331 (85.21)'
$ws.Range("P20").Value = 'maa://37442 (95.12)'
$ws.Range("AF21").Value = 'maa://22524 (94.5), *maa://22432 (77.14)'
$ws.Range("L22").Value = 'maa://27127 (80.87), *maa://22751 (71.01)'
# Leading apostrophe forces this numeric-looking value to be stored as text
# (matches the source cell's existing text type instead of becoming a number).
$ws.Range("C24").Value = '''2'
$ws.Range("D24").Value = '*maa://24368 (78.04), **maa://46650 (50.0)'
$ws.Range("X24").Value = 'maa://29988 (84.74), maa://23504 (93.29), **maa://22892 (40.14), *maa://25141 (76.92), *maa://36663 (78.67), ***maa://22815 (23.08)'
$ws.Range("AF27").Value = 'maa://24023 (97.26)'
$ws.Range("X28").Value = 'maa://39929 (90.5), maa://41749 (90.12), ***maa://39723 (13.89)'
$ws.Range("AF28").Value = 'maa://36660 (92.13), *maa://36701 (65.52)'
$ws.Range("P29").Value = '*maa://23168 (58.06), *maa://30050 (51.52)'
$ws.Range("AB30").Value = 'maa://42979 (96.71), maa://45822 (100.0), *maa://45045 (80.0)'
$ws.Range("T32").Value = 'maa://42859 (96.49), maa://41108 (88.0), maa://41238 (97.09), maa://45523 (100.0)'
$ws.Range("L35").Value = 'maa://41296 (96.13)'
$ws.Range("L37").Value = 'maa://45718 (99.04), maa://45789 (100.0)'
$ws.Range("AF38").Value = 'maa://36697 (86.41)'
$ws.Range("H39").Value = 'maa://36670 (88.66), maa://25199 (84.82), maa://30434 (91.03), ***maa://25036 (16.0), *maa://45059 (80.0), *maa://44165 (66.67)'
$ws.Range("T39").Value = 'maa://45788 (81.4), maa://45790 (81.82)'
$ws.Range("H46").Value = 'maa://35931 (92.56), maa://43901 (90.48)'
$ws.Range("H47").Value = 'maa://27410 (96.42), maa://29661 (97.28), maa://28038 (84.62)'
$ws.Range("H55").Value = 'maa://32532 (92.15)'
$ws.Range("H59").Value = 'maa://27746 (82.3), maa://31270 (95.28)'
